$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column K (11), shifting Sample ID, etc. one column right
$ws.Columns.Item(11).Insert()

# Rename existing J1 header and set the new K1 header
# (set K1 first so new shared strings are appended in "exposed_height",
# "submerged_depth" order, matching the target workbook)
$ws.Range("K1").Value = "exposed_height"
$ws.Range("J1").Value = "submerged_depth"

# Match the author's final view/selection state
$ws.Range("K2").Select()
$excel.ActiveWindow.ScrollColumn = 5
